$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9933785796165466
$ws.Range("B1").Value = 1.78658127784729
$ws.Range("C1").Value = 4.999735355377197
$ws.Range("D1").Value = 1.376420736312866
$ws.Range("E1").Value = 1.320993423461914
